$d = $word.ActiveDocument

# 1. Update the letter date
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing-address line "969 Story Road, San Jose CA 95122" into
#    two separate paragraphs: "969 Story Road" and "San Jose, CA 95122".
#    Only the first occurrence (the mailing address block) should be split;
#    the "PROPERTY ADDRESS" table cell further down keeps the original text.
$rng = $d.Content
$rng.Find.Execute("969 Story Road, San Jose CA 95122", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null
$rng.Text = "969 Story Road`rSan Jose, CA 95122"

# 3. Remove the now-empty "NoSpacing" paragraph that directly follows the
#    "Board of Directors" line.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "`r" -and $p.Style.NameLocal -eq "No Spacing") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text -like "*Board of Directors*") {
            $p.Range.Delete()
            break
        }
    }
}
